$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Fix typo in C51 (missing leading "V" in "Verify")
$ws.Range("C51").Value = "Verify that Help file is accessible from within the application via the 'Help' link on the profile menu||Verify that system makes available the WAT help file (PDF) as a static link||Verify that user should able to access feedback survey page using 'Feedback' link on the profile menu"

# Bring formatting for the two new rows (53 & 54) in line with the existing
# table rows by copying the row immediately above them.
$ws.Range("A52:E52").Copy($ws.Range("A53:E53"))
$ws.Range("A52:E52").Copy($ws.Range("A54:E54"))

# Column D on row 53 uses the "vertical-top" wrapped style (same as D47/D48/...)
$ws.Range("D47").Copy($ws.Range("D53"))

# Row 53 - WAT55
$ws.Range("A53").Value = "WAT55"
$ws.Range("B53").Value = "WAT-730"
$ws.Range("C53").Value = "Verify that ‘Account’ settings link opens up the account settings modal."
$ws.Range("D53").Value = "Y"

# Row 54 - WAT56
$ws.Range("A54").Value = "WAT56"
$ws.Range("B54").Value = "WAT-733"
$ws.Range("C54").Value = "Verify that ‘Terms of Use’ link takes you to /#/terms-of-use-app."
$ws.Range("D54").Value = "Y"

# Keep dimension / view roughly in sync with the new data extent.
$ws.Range("C54").Select()
